# Update the "through December N" cutoff date from Dec 14 to Dec 15
# (workbook tracks the current running month, with one comparison column
# per year, each cut off at the same day-of-month so totals are comparable).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "Through 2022-12-15"

# Update the header label for the current (December 2022) column.
$ws.Range("B1").Value = "December 2022 (through December 15)"

# Update the per-neighborhood counts for every "December" column
# (B = Dec 2022, N = Dec 2021, Z = Dec 2020, AL = Dec 2019,
#  AX = Dec 2018, BJ = Dec 2017, BV = Dec 2016, CH = Dec 2015)
# to reflect one additional day (Dec 15) of data for each year.

$changes = @{
    "AL2"  = 1
    "BJ2"  = 7
    "Z3"   = 2
    "N4"   = 5
    "N5"   = 2
    "Z5"   = 1
    "Z6"   = 1
    "N7"   = 2
    "AX8"  = 2
    "N9"   = 4
    "Z9"   = 5
    "Z10"  = 3
    "BJ10" = 1
    "N14"  = 4
    "AX14" = 7
    "BJ14" = 5
    "BV14" = 4
    "CH14" = 2
    "B20"  = 3
    "AL20" = 2
    "BJ20" = 3
    "N26"  = 3
    "B41"  = 3
    "N41"  = 3
    "Z42"  = 2
    "BJ45" = 2
    "BV48" = 2
    "N55"  = 2
    "N67"  = 2
    "B76"  = 1
    "N92"  = 2
    "Z92"  = 2
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
